# Correções turma A e C.
# Insert a new column "Resenha Novos Keynesianos" before the existing
# "Email" column (shifts Email from column J to column K) and fill in
# the new column's header + per-student scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at J; everything from J onward (incl. "Email")
# shifts one column to the right (J -> K).
$ws.Columns("J").Insert()

# New header for the inserted column.
$ws.Range("J1").Value = "Resenha Novos Keynesianos"

# Per-student "Resenha Novos Keynesianos" scores (row number -> value).
# Rows not present here (5, 14, 15, 16) are left blank, matching students
# who have no grade recorded for this assignment.
$scores = @{
    2  = 5
    3  = 7
    4  = 0
    6  = 0
    7  = 7
    8  = 7
    9  = 7
    10 = 7
    11 = 5
    12 = 7
    13 = 0
    17 = 10
    18 = 3
    19 = 10
    20 = 5
    21 = 0
    22 = 7
    23 = 5
    24 = 10
    25 = 7
    26 = 0
    27 = 10
    28 = 7
    29 = 0
    30 = 7
    31 = 0
    32 = 0
    33 = 5
    34 = 0
    35 = 0
    36 = 7
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 3
    42 = 0
    43 = 7
    44 = 0
}

foreach ($row in $scores.Keys) {
    $ws.Range("J$row").Value = $scores[$row]
}
